$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "validation base with gap to select features"
$ws.Range("B2").Value = 2.02

$ws.Range("A2").Select()
